$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "40÷3=13, 1"
$t.Cell(1, 2).Range.Text = "96÷6=16, 0"
$t.Cell(1, 3).Range.Text = "14÷9=1, 5"
$t.Cell(1, 4).Range.Text = "38÷7=5, 3"
$t.Cell(1, 5).Range.Text = "81÷2=40, 1"
$t.Cell(5, 1).Range.Text = "20÷4=5, 0"
$t.Cell(5, 2).Range.Text = "95÷4=23, 3"
$t.Cell(5, 3).Range.Text = "40÷2=20, 0"
$t.Cell(5, 4).Range.Text = "78÷3=26, 0"
$t.Cell(5, 5).Range.Text = "54÷9=6, 0"
$t.Cell(9, 1).Range.Text = "84÷9=9, 3"
$t.Cell(9, 2).Range.Text = "59÷2=29, 1"
$t.Cell(9, 3).Range.Text = "23÷8=2, 7"
$t.Cell(9, 4).Range.Text = "71÷4=17, 3"
$t.Cell(9, 5).Range.Text = "96÷9=10, 6"
$t.Cell(13, 1).Range.Text = "90÷9=10, 0"
$t.Cell(13, 2).Range.Text = "77÷4=19, 1"
$t.Cell(13, 3).Range.Text = "81÷7=11, 4"
$t.Cell(13, 4).Range.Text = "46÷4=11, 2"
$t.Cell(13, 5).Range.Text = "85÷8=10, 5"
$t.Cell(17, 1).Range.Text = "98÷8=12, 2"
$t.Cell(17, 2).Range.Text = "11÷3=3, 2"
$t.Cell(17, 3).Range.Text = "77÷6=12, 5"
$t.Cell(17, 4).Range.Text = "65÷5=13, 0"
$t.Cell(17, 5).Range.Text = "10÷3=3, 1"
